$p = $ppt.ActivePresentation

# Slide 1: Title slide - update main title text (first run of the title text box)
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item("Title 1")
$titleShape.TextFrame.TextRange.Runs(1,1).Text = "Bayesian Analysis of PGA Golf Scores with Stan"

# Slide 2: Agenda slide - update bullet text
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item("Content Placeholder 2")
$tf = $contentShape.TextFrame
$tr = $tf.TextRange

$tr.Paragraphs(2,1).Runs(1,1).Text = "Multi-level Model of Scores and Tournaments"
$tr.Paragraphs(3,1).Runs(1,1).Text = "Stokes Gained Intro and Score Regression"
$tr.Paragraphs(4,1).Runs(1,1).Text = "Time Series Modeling of Strokes Gained "
